# Generate Report for Handoff
#
# The localization-status report is regenerated with a new handoff id
# (67df7530-5b68-4b51-a5a2-008c033e794d -> b9cd0c30-cc44-4e04-86d7-2dd77cac3972),
# a new target-file hash (e7712af5ed7369c43716f7b10bb5f16d507fb0af ->
# c328a061229b321b6cf352ef241b3f06dc797047), and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = "b9cd0c30-cc44-4e04-86d7-2dd77cac3972.md"
$overview.Range("D2").Value = "2016-47-18 14:47:59"

# --- "zh-cn" sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = "b9cd0c30-cc44-4e04-86d7-2dd77cac3972.md"
$zhcn.Range("D2").Value = "b9cd0c30-cc44-4e04-86d7-2dd77cac3972.c328a061229b321b6cf352ef241b3f06dc797047.zh-cn.xlf"
$zhcn.Range("E2").Value = "2016-03-18 14:47:56"

# --- "de-de" sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = "b9cd0c30-cc44-4e04-86d7-2dd77cac3972.md"
$dede.Range("D2").Value = "b9cd0c30-cc44-4e04-86d7-2dd77cac3972.c328a061229b321b6cf352ef241b3f06dc797047.de-de.xlf"
$dede.Range("E2").Value = "2016-03-18 14:47:59"
